# Apply the "Office Theme" design/colour scheme to the deck's (sole
# reachable) theme part, ppt/theme/theme1.xml, which is used by the
# single Slide Master (and therefore every slide layout / slide).
#
# This mirrors picking a new theme from PowerPoint's Design tab: the
# 12 theme colours (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) change
# from the previous "Integral" (Red Violet) palette to the standard
# Office palette. The font scheme and format scheme (fills/lines/
# effects) are already identical between the old and new themes, so
# no other theme element needs touching.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

# Office theme colours, converted from RRGGBB hex to the BGR-ordered
# long integer PowerPoint's RGB property expects.
$colorScheme.Item(1).RGB  = 0          # dk1      000000
$colorScheme.Item(2).RGB  = 16777215   # lt1      FFFFFF
$colorScheme.Item(3).RGB  = 6968388    # dk2      44546A
$colorScheme.Item(4).RGB  = 15132391   # lt2      E7E6E6
$colorScheme.Item(5).RGB  = 13998939   # accent1  5B9BD5
$colorScheme.Item(6).RGB  = 3243501    # accent2  ED7D31
$colorScheme.Item(7).RGB  = 10855845   # accent3  A5A5A5
$colorScheme.Item(8).RGB  = 49407      # accent4  FFC000
$colorScheme.Item(9).RGB  = 12874308   # accent5  4472C4
$colorScheme.Item(10).RGB = 4697456    # accent6  70AD47
$colorScheme.Item(11).RGB = 12673797   # hlink    0563C1
$colorScheme.Item(12).RGB = 7491477    # folHlink 954F72
